# Fix the 2050 column header label (was showing a stray numeric value left
# over from a copy/paste) and remove the "Total" rows from the scenario
# tables.

$wb = $excel.ActiveWorkbook

# Sheets 1-4: last-year header (column E, row 1) and trailing "Total" row
# (row 13) need fixing/removal.
$sheetsWithTotalRow = @(1, 2, 3, 4)

foreach ($idx in $sheetsWithTotalRow) {
    $ws = $wb.Worksheets.Item($idx)

    # Fix the label in E1 based on what D1 currently contains (the other
    # sheets use single years, the "Incremental" sheet uses year ranges).
    $d1 = $ws.Cells.Item(1, 4).Value
    if ($d1 -eq "2031-2040") {
        $ws.Cells.Item(1, 5).Value = "'2041-2050"
    } else {
        $ws.Cells.Item(1, 5).Value = "'2050"
    }

    # Remove the "Total" row (row 13).
    $ws.Rows.Item(13).Delete()
}

# Sheet 5: only the E1 label needs fixing, it has no Total row.
$ws5 = $wb.Worksheets.Item(5)
$ws5.Cells.Item(1, 5).Value = "'2050"

# Sheet 6: only the trailing "Total" row (row 4) needs removal.
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
